# report.xlsx update: "Component" -> "Source" column, housekeeping on sheet layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")
$ws2 = $wb.Worksheets.Item("Sheet3")

# --- Header / template text changes -----------------------------------
# Header row (row 3): "Component" -> "Source"
$ws.Cells.Item(3, 3).Value = "Source"
# Template row (row 4): "${issue.component}" -> "${issue.source}"
$ws.Cells.Item(4, 3).Value = '${issue.source}'

# --- Remove the stray trailing empty row (was row 1048576) -------------
$ws.Rows.Item(1048576).Delete()

# --- Row heights ---------------------------------------------------------
# Data row height now matches the new default (13.8pt)
$ws.Rows.Item(4).RowHeight = 13.8

# --- Column widths ---------------------------------------------------------
# Columns A and B become one uniform-width block, C and D are resized.
$ws.Range("A1:B1").EntireColumn.ColumnWidth = 17.333333333333332
$ws.Columns.Item(3).ColumnWidth = 47
$ws.Columns.Item(4).ColumnWidth = 57.5

$ws2.Columns.Item(1).ColumnWidth = 7.5

# --- Selection -------------------------------------------------------------
$ws.Range("D20").Select()

Write-Host "C3:" $ws.Cells.Item(3, 3).Value2
Write-Host "C4:" $ws.Cells.Item(4, 3).Value2
